# === Nexial "#system" sheet: add storeKeys(json,jsonpath,var) to the json
# function list (shared with the `keys(jsonpath)` JSON expression), and
# remove the stray one-cell "text" column, shifting web/webalert/webcookie/
# ws/ws.async/xml one column to the left. ===

$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

# Clear every column whose data is being rewritten/reshuffled below.
$clearRanges = @(
    "A1:A31",
    "M1:M18",
    "Y1:Y129",
    "Z1:Z129",
    "AA1:AA8",
    "AB1:AB17",
    "AC1:AC17",
    "AD1:AD27",
    "AE1:AE27"
)
foreach ($r in $clearRanges) {
    $sys.Range($r).ClearContents()
}

# Final resolved values for the affected cells (address, value).
$cellValues = @(
    @("A1", "target"),
    @("A2", "aws.s3"),
    @("A3", "aws.ses"),
    @("A4", "aws.sqs"),
    @("A5", "base"),
    @("A6", "csv"),
    @("A7", "desktop"),
    @("A8", "excel"),
    @("A9", "external"),
    @("A10", "image"),
    @("A11", "io"),
    @("A12", "jms"),
    @("A13", "json"),
    @("A14", "localdb"),
    @("A15", "macro"),
    @("A16", "mail"),
    @("A17", "number"),
    @("A18", "pdf"),
    @("A19", "rdbms"),
    @("A20", "redis"),
    @("A21", "sms"),
    @("A22", "sound"),
    @("A23", "ssh"),
    @("A24", "step"),
    @("A25", "web"),
    @("A26", "webalert"),
    @("A27", "webcookie"),
    @("A28", "ws"),
    @("A29", "ws.async"),
    @("A30", "xml"),
    @("M1", "json"),
    @("M2", "addOrReplace(json,jsonpath,input,var)"),
    @("M3", "assertCorrectness(json,schema)"),
    @("M4", "assertElementCount(json,jsonpath,count)"),
    @("M5", "assertElementNotPresent(json,jsonpath)"),
    @("M6", "assertElementPresent(json,jsonpath)"),
    @("M7", "assertEqual(expected,actual)"),
    @("M8", "assertValue(json,jsonpath,expected)"),
    @("M9", "assertValues(json,jsonpath,array,exactOrder)"),
    @("M10", "assertWellformed(json)"),
    @("M11", "beautify(json,var)"),
    @("M12", "compact(var,json,removeEmpty)"),
    @("M13", "fromCsv(csv,header,jsonFile)"),
    @("M14", "minify(json,var)"),
    @("M15", "storeCount(json,jsonpath,var)"),
    @("M16", "storeKeys(json,jsonpath,var)"),
    @("M17", "storeValue(json,jsonpath,var)"),
    @("M18", "storeValues(json,jsonpath,var)"),
    @("Y1", "web"),
    @("Y2", "assertAndClick(locator,label)"),
    @("Y3", "assertAttribute(locator,attrName,value)"),
    @("Y4", "assertAttributeContains(locator,attrName,contains)"),
    @("Y5", "assertAttributeNotContains(locator,attrName,contains)"),
    @("Y6", "assertAttributeNotPresent(locator,attrName)"),
    @("Y7", "assertAttributePresent(locator,attrName)"),
    @("Y8", "assertChecked(locator)"),
    @("Y9", "assertContainCount(locator,text,count)"),
    @("Y10", "assertCssNotPresent(locator,property)"),
    @("Y11", "assertCssPresent(locator,property,value)"),
    @("Y12", "assertElementByAttributes(nameValues)"),
    @("Y13", "assertElementByText(locator,text)"),
    @("Y14", "assertElementCount(locator,count)"),
    @("Y15", "assertElementNotPresent(locator)"),
    @("Y16", "assertElementPresent(locator)"),
    @("Y17", "assertElementsPresent(prefix)"),
    @("Y18", "assertFocus(locator)"),
    @("Y19", "assertFrameCount(count)"),
    @("Y20", "assertFramePresent(frameName)"),
    @("Y21", "assertIECompatMode()"),
    @("Y22", "assertIENativeMode()"),
    @("Y23", "assertLinkByLabel(label)"),
    @("Y24", "assertNotChecked(locator)"),
    @("Y25", "assertNotFocus(locator)"),
    @("Y26", "assertNotText(locator,text)"),
    @("Y27", "assertNotVisible(locator)"),
    @("Y28", "assertOneMatch(locator)"),
    @("Y29", "assertScrollbarHNotPresent(locator)"),
    @("Y30", "assertScrollbarHPresent(locator)"),
    @("Y31", "assertScrollbarVNotPresent(locator)"),
    @("Y32", "assertScrollbarVPresent(locator)"),
    @("Y33", "assertTable(locator,row,column,text)"),
    @("Y34", "assertText(locator,text)"),
    @("Y35", "assertTextContains(locator,text)"),
    @("Y36", "assertTextCount(locator,text,count)"),
    @("Y37", "assertTextList(locator,list,ignoreOrder)"),
    @("Y38", "assertTextMatches(text,minMatch,scrollTo)"),
    @("Y39", "assertTextNotContains(locator,text)"),
    @("Y40", "assertTextNotPresent(text)"),
    @("Y41", "assertTextOrder(locator,descending)"),
    @("Y42", "assertTextPresent(text)"),
    @("Y43", "assertTitle(text)"),
    @("Y44", "assertValue(locator,value)"),
    @("Y45", "assertValueOrder(locator,descending)"),
    @("Y46", "assertVisible(locator)"),
    @("Y47", "checkAll(locator)"),
    @("Y48", "clearLocalStorage()"),
    @("Y49", "click(locator)"),
    @("Y50", "clickAll(locator)"),
    @("Y51", "clickAndWait(locator,waitMs)"),
    @("Y52", "clickByLabel(label)"),
    @("Y53", "clickByLabelAndWait(label,waitMs)"),
    @("Y54", "clickOffset(locator,x,y)"),
    @("Y55", "clickWithKeys(locator,keys)"),
    @("Y56", "close()"),
    @("Y57", "closeAll()"),
    @("Y58", "deselect(locator,text)"),
    @("Y59", "deselectMulti(locator,array)"),
    @("Y60", "dismissInvalidCert()"),
    @("Y61", "dismissInvalidCertPopup()"),
    @("Y62", "doubleClick(locator)"),
    @("Y63", "doubleClickAndWait(locator,waitMs)"),
    @("Y64", "doubleClickByLabel(label)"),
    @("Y65", "doubleClickByLabelAndWait(label,waitMs)"),
    @("Y66", "dragAndDrop(fromLocator,toLocator)"),
    @("Y67", "dragTo(fromLocator,xOffset,yOffset)"),
    @("Y68", "editLocalStorage(key,value)"),
    @("Y69", "executeScript(var,script)"),
    @("Y70", "focus(locator)"),
    @("Y71", "goBack()"),
    @("Y72", "goBackAndWait()"),
    @("Y73", "maximizeWindow()"),
    @("Y74", "mouseOver(locator)"),
    @("Y75", "open(url)"),
    @("Y76", "openAndWait(url,waitMs)"),
    @("Y77", "openHttpBasic(url,username,password)"),
    @("Y78", "openIgnoreTimeout(url)"),
    @("Y79", "refresh()"),
    @("Y80", "refreshAndWait()"),
    @("Y81", "resizeWindow(width,height)"),
    @("Y82", "rightClick(locator)"),
    @("Y83", "saveAllWindowIds(var)"),
    @("Y84", "saveAllWindowNames(var)"),
    @("Y85", "saveAttribute(var,locator,attrName)"),
    @("Y86", "saveAttributeList(var,locator,attrName)"),
    @("Y87", "saveCount(var,locator)"),
    @("Y88", "saveDivsAsCsv(headers,rows,cells,nextPage,file)"),
    @("Y89", "saveElement(var,locator)"),
    @("Y90", "saveElements(var,locator)"),
    @("Y91", "saveLocalStorage(var,key)"),
    @("Y92", "saveLocation(var)"),
    @("Y93", "savePageAs(var,sessionIdName,url)"),
    @("Y94", "savePageAsFile(sessionIdName,url,file)"),
    @("Y95", "saveTableAsCsv(locator,nextPageLocator,file)"),
    @("Y96", "saveText(var,locator)"),
    @("Y97", "saveTextArray(var,locator)"),
    @("Y98", "saveTextSubstringAfter(var,locator,delim)"),
    @("Y99", "saveTextSubstringBefore(var,locator,delim)"),
    @("Y100", "saveTextSubstringBetween(var,locator,start,end)"),
    @("Y101", "saveValue(var,locator)"),
    @("Y102", "saveValues(var,locator)"),
    @("Y103", "scrollElement(locator,xOffset,yOffset)"),
    @("Y104", "scrollLeft(locator,pixel)"),
    @("Y105", "scrollPage(xOffset,yOffset)"),
    @("Y106", "scrollRight(locator,pixel)"),
    @("Y107", "scrollTo(locator)"),
    @("Y108", "select(locator,text)"),
    @("Y109", "selectFrame(locator)"),
    @("Y110", "selectMulti(locator,array)"),
    @("Y111", "selectMultiOptions(locator)"),
    @("Y112", "selectText(locator)"),
    @("Y113", "selectWindow(winId)"),
    @("Y114", "selectWindowAndWait(winId,waitMs)"),
    @("Y115", "selectWindowByIndex(index)"),
    @("Y116", "selectWindowByIndexAndWait(index,waitMs)"),
    @("Y117", "toggleSelections(locator)"),
    @("Y118", "type(locator,value)"),
    @("Y119", "typeKeys(locator,value)"),
    @("Y120", "uncheckAll(locator)"),
    @("Y121", "unselectAllText()"),
    @("Y122", "upload(fieldLocator,file)"),
    @("Y123", "verifyContainText(locator,text)"),
    @("Y124", "verifyText(locator,text)"),
    @("Y125", "wait(waitMs)"),
    @("Y126", "waitForElementPresent(locator)"),
    @("Y127", "waitForPopUp(winId,waitMs)"),
    @("Y128", "waitForTextPresent(text)"),
    @("Y129", "waitForTitle(text)"),
    @("Z1", "webalert"),
    @("Z2", "accept()"),
    @("Z3", "assertPresent()"),
    @("Z4", "assertText(text,matchBy)"),
    @("Z5", "dismiss()"),
    @("Z6", "replyCancel(text)"),
    @("Z7", "replyOK(text)"),
    @("Z8", "storeText(var)"),
    @("AA1", "webcookie"),
    @("AA2", "assertNotPresent(name)"),
    @("AA3", "assertPresent(name)"),
    @("AA4", "assertValue(name,value)"),
    @("AA5", "delete(name)"),
    @("AA6", "deleteAll()"),
    @("AA7", "save(var,name)"),
    @("AA8", "saveAll(var)"),
    @("AB1", "ws"),
    @("AB2", "assertReturnCode(var,returnCode)"),
    @("AB3", "delete(url,body,var)"),
    @("AB4", "download(url,queryString,saveTo)"),
    @("AB5", "get(url,queryString,var)"),
    @("AB6", "head(url,var)"),
    @("AB7", "header(name,value)"),
    @("AB8", "headerByVar(name,var)"),
    @("AB9", "jwtParse(var,token,key)"),
    @("AB10", "jwtSignHS256(var,payload,key)"),
    @("AB11", "oauth(var,url,auth)"),
    @("AB12", "patch(url,body,var)"),
    @("AB13", "post(url,body,var)"),
    @("AB14", "put(url,body,var)"),
    @("AB15", "saveResponsePayload(var,file,append)"),
    @("AB16", "soap(action,url,payload,var)"),
    @("AB17", "upload(url,body,fileParams,var)"),
    @("AC1", "ws.async"),
    @("AC2", "delete(url,body,output)"),
    @("AC3", "download(url,queryString,saveTo)"),
    @("AC4", "get(url,queryString,output)"),
    @("AC5", "head(url,output)"),
    @("AC6", "patch(url,body,output)"),
    @("AC7", "post(url,body,output)"),
    @("AC8", "put(url,body,output)"),
    @("AD1", "xml"),
    @("AD2", "append(xml,xpath,content,var)"),
    @("AD3", "assertCorrectness(xml,schema)"),
    @("AD4", "assertElementCount(xml,xpath,count)"),
    @("AD5", "assertElementNotPresent(xml,xpath)"),
    @("AD6", "assertElementPresent(xml,xpath)"),
    @("AD7", "assertSoap(wsdl,xml)"),
    @("AD8", "assertSoapFaultCode(expected,xml)"),
    @("AD9", "assertSoapFaultString(expected,xml)"),
    @("AD10", "assertValue(xml,xpath,expected)"),
    @("AD11", "assertValues(xml,xpath,array,exactOrder)"),
    @("AD12", "assertWellformed(xml)"),
    @("AD13", "beautify(xml,var)"),
    @("AD14", "clear(xml,xpath,var)"),
    @("AD15", "delete(xml,xpath,var)"),
    @("AD16", "insertAfter(xml,xpath,content,var)"),
    @("AD17", "insertBefore(xml,xpath,content,var)"),
    @("AD18", "minify(xml,var)"),
    @("AD19", "prepend(xml,xpath,content,var)"),
    @("AD20", "replace(xml,xpath,content,var)"),
    @("AD21", "replaceIn(xml,xpath,content,var)"),
    @("AD22", "storeCount(xml,xpath,var)"),
    @("AD23", "storeSoapFaultCode(var,xml)"),
    @("AD24", "storeSoapFaultDetail(var,xml)"),
    @("AD25", "storeSoapFaultString(var,xml)"),
    @("AD26", "storeValue(xml,xpath,var)"),
    @("AD27", "storeValues(xml,xpath,var)")
)
foreach ($pair in $cellValues) {
    $sys.Range($pair[0]).Value = $pair[1]
}

# Keep the named ranges in the workbook in sync with the new layout.
$newRefersTo = @{
    "json"      = "='#system'!`$M`$2:`$M`$18"
    "target"    = "='#system'!`$A`$2:`$A`$30"
    "web"       = "='#system'!`$Y`$2:`$Y`$129"
    "webalert"  = "='#system'!`$Z`$2:`$Z`$8"
    "webcookie" = "='#system'!`$AA`$2:`$AA`$8"
    "ws"        = "='#system'!`$AB`$2:`$AB`$17"
    "ws.async"  = "='#system'!`$AC`$2:`$AC`$8"
    "xml"       = "='#system'!`$AD`$2:`$AD`$27"
}
foreach ($n in $wb.Names) {
    if ($newRefersTo.ContainsKey($n.Name)) {
        $n.RefersTo = $newRefersTo[$n.Name]
    }
}

